$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The two "is_locked" / "is_enabled" select-list header columns (D, E) are
# removed; the two trailing columns (order_by, rem) shift left to take
# their place, and the now-empty trailing columns (F, G) are deleted.
$orderByVal = $ws.Range("F1").Value()
$remVal = $ws.Range("G1").Value()

$ws.Range("D1").Value = $orderByVal
$ws.Range("E1").Value = $remVal

$ws.Range("F1:G1").Delete(-4159)
